$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a brand new product row at row 21 ("بلاستر مترسيلك 2.5 سم"),
#    pushing the existing rows 21-28 (data rows + totals row) down by one.
# ---------------------------------------------------------------------------
$ws.Rows.Item(21).Insert()

# Re-apply the formatting of the (now shifted) former row 21 - which landed
# on row 22 - onto the freshly inserted, blank row 21 so styles/number
# formats match the rest of the table exactly.
$ws.Range("A22:Q22").Copy()
$ws.Range("A21:Q21").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(21).RowHeight = $ws.Rows.Item(22).RowHeight
$excel.CutCopyMode = 0

# Re-create the merged cells for the new row 21 (mirroring every other data
# row in the table: A:B, C:G, H:K, L:M, N:O).
$ws.Range("A21:B21").Merge()
$ws.Range("C21:G21").Merge()
$ws.Range("H21:K21").Merge()
$ws.Range("L21:M21").Merge()
$ws.Range("N21:O21").Merge()

# Helper scratch cell (kept text-formatted) used so that numeric-looking
# strings ("25.0000", "19:0", ...) are written as plain text, matching the
# shared-string cells already used throughout this sheet, regardless of the
# destination cell's own (numeric) number format.
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"

function Set-TextValue($rangeAddress, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($rangeAddress).PasteSpecial(-4163)  # xlPasteValues
}

# Fill in the new row 21 values.
$ws.Range("A21").Value = 15
Set-TextValue "C21" "بلاستر مترسيلك 2.5 سم"
Set-TextValue "H21" "19:0"
Set-TextValue "L21" "0"
Set-TextValue "N21" "25.00"
Set-TextValue "P21" "25.0000"
Set-TextValue "Q21" "1:0"

# Renumber the (shifted) rows 22-28 so the running index keeps counting
# 16, 17, 18 ... 22 (it used to read 15, 16, 17 ... 21 before the insert).
for ($r = 22; $r -le 28; $r++) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r, 1).Value2 + 1
}

# ---------------------------------------------------------------------------
# 2) Append one more data row (new row 28) before the totals row, repeating
#    what used to be the last product row ("محلول جلوكوز 5%") with the next
#    running index (22).
# ---------------------------------------------------------------------------
$ws.Rows.Item(29).Insert()

$ws.Range("A28:Q28").Copy()
$ws.Range("A29:Q29").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(29).RowHeight = $ws.Rows.Item(28).RowHeight
$excel.CutCopyMode = 0

$ws.Range("A29:B29").Merge()
$ws.Range("C29:G29").Merge()
$ws.Range("H29:K29").Merge()
$ws.Range("L29:M29").Merge()
$ws.Range("N29:O29").Merge()

$ws.Range("A29").Value = 22
Set-TextValue "C29" "محلول جلوكوز 5%"
Set-TextValue "H29" "20:0"
Set-TextValue "L29" "0"
Set-TextValue "N29" "27.00"
Set-TextValue "P29" "27.0000"
Set-TextValue "Q29" "1:0"

$scratch.Clear()

# ---------------------------------------------------------------------------
# 3) Update the grand-total (now on row 30) to reflect the extra row, and
#    refresh the generated timestamp in the footer (now on row 31).
# ---------------------------------------------------------------------------
$ws.Range("P30").Value = 1054.9200000000001

$ws.Range("A31").Value = "Friday, 12 September, 2025 5:31 PM"
